$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = [double]"-5.553054260865053e-08"
$ws.Range("D2").Value = [double]"-1.094440884230607e-05"
$ws.Range("F2").Value = [double]"-5.553054260865053e-08"
$ws.Range("G2").Value = [double]"-2.776527132652973e-07"
$ws.Range("C3").Value = [double]"2.848439828451177e-08"
$ws.Range("D3").Value = [double]"1.139375931380471e-07"
$ws.Range("E3").Value = [double]"3.482706989466067e-08"
$ws.Range("F3").Value = [double]"5.279848469541548e-08"
$ws.Range("G3").Value = [double]"1.424219915335811e-07"
$ws.Range("C4").Value = [double]"1.976225101098239e-09"
$ws.Range("D4").Value = [double]"7.904900404392956e-09"
$ws.Range("E4").Value = [double]"9.038725401921965e-06"
$ws.Range("F4").Value = [double]"-0.04934920092854345"
$ws.Range("G4").Value = [double]"9.88112550809328e-09"
$ws.Range("C5").Value = [double]"-0.0001235393247043248"
$ws.Range("D5").Value = [double]"-2.074978473487477e-07"
$ws.Range("F5").Value = [double]"-5.187446183718691e-08"
$ws.Range("G5").Value = [double]"-2.593723094079792e-07"
$ws.Range("C6").Value = [double]"-5.457296359256958e-08"
$ws.Range("D6").Value = [double]"-1.075568139796701e-05"
$ws.Range("F6").Value = [double]"-5.457296359256958e-08"
$ws.Range("G6").Value = [double]"-2.728648179628479e-07"
$ws.Range("C7").Value = [double]"-7.270126274605548e-08"
$ws.Range("D7").Value = [double]"-2.908050509842219e-07"
$ws.Range("E7").Value = [double]"-0.000332516144226247"
$ws.Range("F7").Value = [double]"-5.37721321478557e-05"
$ws.Range("G7").Value = [double]"-3.635063137163996e-07"
$ws.Range("B8").Value = [double]"-5.797333824375528e-05"
$ws.Range("C8").Value = [double]"-3.190793664842317e-07"
$ws.Range("D8").Value = [double]"-3.109940143986023e-05"
$ws.Range("F8").Value = [double]"-3.190793664842317e-07"
$ws.Range("G8").Value = [double]"-1.59539683153298e-06"
$ws.Range("C9").Value = [double]"-2.310664770632798e-08"
$ws.Range("D9").Value = [double]"-1.991017529689998e-05"
$ws.Range("F9").Value = [double]"-2.310664770632798e-08"
$ws.Range("G9").Value = [double]"-1.155332380875507e-07"
$ws.Range("C10").Value = [double]"-5.504249038779108e-07"
$ws.Range("D10").Value = [double]"-7.422426165248908e-06"
$ws.Range("E10").Value = [double]"-1.333445055706761e-06"
$ws.Range("F10").Value = [double]"-1.702503624301244e-06"
$ws.Range("G10").Value = [double]"-2.75212451583684e-06"
$ws.Range("C11").Value = [double]"-9.195836803144175e-08"
$ws.Range("D11").Value = [double]"-1.065187200310902e-05"
$ws.Range("F11").Value = [double]"-5.909425453864969e-06"
$ws.Range("G11").Value = [double]"-7.067260685289511e-05"
